$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final reported counts for range B2:J46 (columns B..J, rows 2..46),
# reflecting the finalized report data delivered to the user.
$data = @(
    @(44,32,17,3,2,7,30,3,20),
    @(11,9,6,0,1,3,12,0,4),
    @(1,0,0,0,0,1,3,0,2),
    @(2,0,1,0,0,0,1,0,0),
    @(0,0,1,0,0,0,2,0,0),
    @(3,2,0,0,0,0,0,0,1),
    @(4,2,1,0,1,1,3,0,0),
    @(0,0,2,0,0,0,3,0,0),
    @(1,3,1,0,0,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,2,0,0,0,0,0,0,1),
    @(11,13,2,0,0,0,5,1,2),
    @(6,5,1,0,0,0,1,1,0),
    @(3,4,1,0,0,0,1,0,1),
    @(2,1,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,3,0,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,2,0,0),
    @(10,1,0,2,0,0,8,1,3),
    @(1,1,0,0,0,0,5,1,3),
    @(2,0,0,0,0,0,2,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(7,0,0,2,0,0,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(4,1,4,0,0,2,3,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(2,1,1,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(2,0,3,0,0,2,3,0,0),
    @(5,5,4,0,1,1,2,1,11),
    @(2,3,1,0,1,0,1,0,0),
    @(1,1,0,0,0,0,0,0,7),
    @(0,0,0,0,0,0,0,0,0),
    @(2,0,0,0,0,1,1,1,0),
    @(0,0,0,0,0,0,0,0,3),
    @(0,1,3,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(3,3,1,1,0,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(3,2,1,1,0,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,1,0,0,0,0,0,0,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $rowValues[$j]
    }
}
